$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

Set-TextValue $ws.Range('D2') '255.08'
Set-TextValue $ws.Range('G2') '8'
Set-TextValue $ws.Range('D3') '28.10'
Set-TextValue $ws.Range('E3') '-5.90%'
Set-TextValue $ws.Range('G3') '8'
Set-TextValue $ws.Range('D4') '5.385'
Set-TextValue $ws.Range('E4') '4.37%'
Set-TextValue $ws.Range('G4') '8'
Set-TextValue $ws.Range('D5') '0.05850'
Set-TextValue $ws.Range('E5') '1.01%'
Set-TextValue $ws.Range('G5') '8'
Set-TextValue $ws.Range('D6') '6.708'
Set-TextValue $ws.Range('E6') '0.91%'
Set-TextValue $ws.Range('G6') '8'
Set-TextValue $ws.Range('D7') '0.8665'
Set-TextValue $ws.Range('E7') '1.96%'
Set-TextValue $ws.Range('G7') '8'
Set-TextValue $ws.Range('D8') '0.9092'
Set-TextValue $ws.Range('E8') '5.60%'
Set-TextValue $ws.Range('G8') '8'
Set-TextValue $ws.Range('D9') '0.1419'
Set-TextValue $ws.Range('E9') '2.57%'
Set-TextValue $ws.Range('G9') '8'
Set-TextValue $ws.Range('D10') '0.07153'
Set-TextValue $ws.Range('E10') '0.68%'
Set-TextValue $ws.Range('G10') '8'
Set-TextValue $ws.Range('D11') '0.03181'
Set-TextValue $ws.Range('E11') '-2.19%'
Set-TextValue $ws.Range('G11') '8'
Set-TextValue $ws.Range('D12') '0.09226'
Set-TextValue $ws.Range('E12') '-1.55%'
Set-TextValue $ws.Range('G12') '8'
Set-TextValue $ws.Range('D13') '0.001537'
Set-TextValue $ws.Range('E13') '-0.06%'
Set-TextValue $ws.Range('G13') '8'
Set-TextValue $ws.Range('D14') '0.0006066'
Set-TextValue $ws.Range('E14') '-94.09%'
Set-TextValue $ws.Range('G14') '8'
Set-TextValue $ws.Range('D15') '0.005804'
Set-TextValue $ws.Range('E15') '-4.25%'
Set-TextValue $ws.Range('G15') '8'
Set-TextValue $ws.Range('D16') '3.498'
Set-TextValue $ws.Range('E16') '-0.17%'
Set-TextValue $ws.Range('G16') '8'
Set-TextValue $ws.Range('E17') '0.00%'
Set-TextValue $ws.Range('G17') '8'
Set-TextValue $ws.Range('D18') '2.202'
Set-TextValue $ws.Range('E18') '-0.41%'
Set-TextValue $ws.Range('G18') '8'
Set-TextValue $ws.Range('D19') '0.3174'
Set-TextValue $ws.Range('E19') '0.03%'
Set-TextValue $ws.Range('G19') '8'
Set-TextValue $ws.Range('E20') '2.33%'
Set-TextValue $ws.Range('G20') '8'
Set-TextValue $ws.Range('D21') '0.1316'
Set-TextValue $ws.Range('E21') '1.50%'
Set-TextValue $ws.Range('G21') '8'
Set-TextValue $ws.Range('D22') '3.521'
Set-TextValue $ws.Range('E22') '1.05%'
Set-TextValue $ws.Range('G22') '8'
Set-TextValue $ws.Range('D23') '0.04150'
Set-TextValue $ws.Range('E23') '0.35%'
Set-TextValue $ws.Range('G23') '8'
Set-TextValue $ws.Range('E24') '-0.18%'
Set-TextValue $ws.Range('G24') '8'
Set-TextValue $ws.Range('D25') '0.005042'
Set-TextValue $ws.Range('E25') '21.76%'
Set-TextValue $ws.Range('G25') '8'
Set-TextValue $ws.Range('D26') '0.001222'
Set-TextValue $ws.Range('E26') '-0.43%'
Set-TextValue $ws.Range('G26') '8'
Set-TextValue $ws.Range('D27') '0.0001199'
Set-TextValue $ws.Range('E27') '-0.07%'
Set-TextValue $ws.Range('G27') '8'
Set-TextValue $ws.Range('D28') '0.0001937'
Set-TextValue $ws.Range('E28') '33.64%'
Set-TextValue $ws.Range('G28') '8'
Set-TextValue $ws.Range('G29') '8'
Set-TextValue $ws.Range('G30') '8'
Set-TextValue $ws.Range('G31') '8'
Set-TextValue $ws.Range('G32') '8'
Set-TextValue $ws.Range('G33') '8'
Set-TextValue $ws.Range('G34') '8'
Set-TextValue $ws.Range('G35') '8'
Set-TextValue $ws.Range('G36') '8'
Set-TextValue $ws.Range('G37') '8'
Set-TextValue $ws.Range('G38') '8'
Set-TextValue $ws.Range('G39') '8'
Set-TextValue $ws.Range('D40') '0.03839'
Set-TextValue $ws.Range('E40') '1.96%'
Set-TextValue $ws.Range('G40') '8'
Set-TextValue $ws.Range('E41') '3.17%'
Set-TextValue $ws.Range('G41') '8'
Set-TextValue $ws.Range('B42') 'CEJI'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws.Range('D42') '0.002337'
Set-TextValue $ws.Range('E42') '6.23%'
Set-TextValue $ws.Range('G42') '8'
Set-TextValue $ws.Range('B43') 'KickToken'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range('D43') '0.002948'
Set-TextValue $ws.Range('E43') '-48.47%'
Set-TextValue $ws.Range('G43') '8'
Set-TextValue $ws.Range('D44') '0.01096'
Set-TextValue $ws.Range('E44') '14.69%'
Set-TextValue $ws.Range('G44') '8'
Set-TextValue $ws.Range('D45') '0.00005218'
Set-TextValue $ws.Range('E45') '-1.97%'
Set-TextValue $ws.Range('G45') '8'
Set-TextValue $ws.Range('E46') '-0.08%'
Set-TextValue $ws.Range('G46') '8'
Set-TextValue $ws.Range('D47') '0.08749'
Set-TextValue $ws.Range('E47') '23.21%'
Set-TextValue $ws.Range('G47') '8'
Set-TextValue $ws.Range('D48') '0.002156'
Set-TextValue $ws.Range('E48') '-1.28%'
Set-TextValue $ws.Range('G48') '8'
Set-TextValue $ws.Range('D49') '0.00002099'
Set-TextValue $ws.Range('E49') '-0.08%'
Set-TextValue $ws.Range('G49') '8'
Set-TextValue $ws.Range('D50') '0.0001999'
Set-TextValue $ws.Range('E50') '-0.08%'
Set-TextValue $ws.Range('G50') '8'
Set-TextValue $ws.Range('G51') '8'
